$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.946.92'
$ws.Range('E2').Value = '  -0.33%  '
$ws.Range('D3').Value = '1.673.48'
$ws.Range('E3').Value = '  +0.87%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.90'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.60%  '
$ws.Range('E6').Value = '  +1.81%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  +0.17%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.26'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.64%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0887'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.59%  '
$ws.Range('D12').Value = '1.909.83'
$ws.Range('E12').Value = '  +0.90%  '
$ws.Range('D13').Value = '1.673.52'
$ws.Range('E13').Value = '  +1.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.08'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.17%  '
$ws.Range('E15').Value = '  +0.57%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.62'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.08%  '
$ws.Range('D17').Value = '26.945.46'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '8.10'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.52%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '235.06'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.83%  '
$ws.Range('D20').Value = '0.0₃0734'
$ws.Range('E20').Value = '  -0.44%  '
$ws.Range('E21').Value = '  +0.17%  '
$ws.Range('E22').Value = '  +0.37%  '
$ws.Range('E23').Value = '  -1.27%  '
$ws.Range('E24').Value = '  -2.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.42'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.18%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.21'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.02'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.88%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.113'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.04%  '
$ws.Range('E29').Value = '  +0.18%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0497'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('E31').Value = '  -0.27%  '
$ws.Range('E32').Value = '  +0.47%  '
$ws.Range('D33').Value = '1.485.97'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.13'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.76%  '
$ws.Range('E35').Value = '  +3.00%  '
$ws.Range('E36').Value = '  +0.32%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.584'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.49%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.895'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.69%  '
$ws.Range('E39').Value = '  +0.45%  '
$ws.Range('E40').Value = '  +8.18%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.85'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E42').Value = '  +0.13%  '
$ws.Range('E43').Value = '  +2.83%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '67.47'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.50%  '
$ws.Range('D45').Value = '1.815.28'
$ws.Range('E45').Value = '  +0.64%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.776'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '90.53'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.21%  '
$ws.Range('E48').Value = '  -0.04%  '
$ws.Range('E49').Value = '  +1.64%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.73'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.42%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0508'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.18%  '
